$d = $word.ActiveDocument

# Step 1: "being" -> "to be"
$rng = $d.Content
$rng.Find.Execute("being", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "to be"
$rng.Font.Name = "Times New Roman"

# Step 2: find the boundary between "submitted to" and "EuroVis", then delete the
# lone space run that sits between them
$r1 = $d.Content
$r1.Find.Execute("submitted to", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterSubmitted = $r1.End

$r2 = $d.Content
$r2.Find.Execute("EuroVis", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$beforeEuroVis = $r2.Start

$spaceRng = $d.Range($afterSubmitted, $beforeEuroVis)
$spaceRng.Text = ""

# Step 3: "EuroVis" -> " IEEE TVCG"
$r3 = $d.Content
$r3.Find.Execute("EuroVis", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.Text = " IEEE TVCG"
$r3.Font.Name = "Times New Roman"

# Step 4: move the _GoBack bookmark so it sits right before the sentence's final period
$r4 = $d.Content
$r4.Find.Execute("IEEE TVCG.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$periodPos = $r4.End - 1
$target = $d.Range($periodPos, $periodPos)
$d.Bookmarks.Add("_GoBack", $target)
